# Hortaliza, Femacal de La Calera - Apio
# A new daily price record is inserted as row 341 (pushing the existing
# rows 341-454 down to 342-455), growing the used range from A1:R454 to
# A1:R455.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 341..454 down by one to make room for the new record.
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new observation.
$ws.Range("A341").Value = 3
$ws.Range("B341").Value = "Femacal de La Calera"
$ws.Range("C341").Value = "Coquimbo"
$ws.Range("D341").Value = 44809
$ws.Range("E341").Value = 5
$ws.Range("F341").Value = 100112017
$ws.Range("G341").Value = "Apio"
$ws.Range("H341").Value = "Americana (o)"
$ws.Range("I341").Value = "Primera"
$ws.Range("J341").Value = 210
$ws.Range("K341").Value = 9500
$ws.Range("L341").Value = 10000
$ws.Range("M341").Value = 9786
$ws.Range("N341").Value = "$/docena de matas"
$ws.Range("O341").Value = "Pan de Azúcar"
$ws.Range("P341").Value = 1631
$ws.Range("Q341").Value = 6
$ws.Range("R341").Value = "Hortaliza"
